$d = $word.ActiveDocument

# Locate the "12.6.2025" learning-diary entry paragraph and update its text.
# Old: "Started the course project. I already set up the API routes and next
#       I'm going to set up the database with MongoDB."
# New: "Started the course project. Goals Setter Application. I already set
#       up the API routes and the REST API for Goals. Tomorrow will do the
#       Users."

$apos = [char]8217  # curly right single quote used in the original text

# 1) Insert the new "Goals Setter Application." sentence right after the
#    opening sentence.
$d.Content.Find.Execute(
    "Started the course project. I already ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Started the course project. Goals Setter Application. I already ",
    2) | Out-Null

# 2) Replace the remainder of the paragraph with the updated continuation.
$d.Content.Find.Execute(
    "set up the API routes and next I" + $apos + "m going to set up the database with MongoDB.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "set up the API routes and the REST API for Goals. Tomorrow will do the Users.",
    2) | Out-Null
